$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.268.69"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.847.18"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6736"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07459"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2954"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.99"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07717"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.849.85"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.012"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6739"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.16"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.170"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "29.264.73"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.206"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "161.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.728"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1412"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.05"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.515"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05336"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7595"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.875"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.678"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "1.322.64"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01804"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.728"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9237"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08284"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +10.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.54"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").Value = "1.990.13"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.779"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.193"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05963"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.32%  "
